# 开发进度.xlsx — "add element EValueStick to QStudioSCADA"
#
# The 画面编辑器 (Screen Editor) progress sheet tracks development status for
# each UI element in column B ("状态" / Status), using a dropdown of
# 未开始 (Not Started, red) / 进行中 (In Progress, yellow) / 已完成 (Completed, green).
#
# EValueStick == 数值棒图 (row 14): work on it is now finished, so its status
# moves from "进行中" (In Progress) to "已完成" (Completed).
# Work then begins on the next element, 罐形容器 (row 15, Tank container),
# whose status moves from "未开始" (Not Started) to "进行中" (In Progress).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("画面编辑器")

# Row 14 - 数值棒图 (EValueStick): In Progress -> Completed.
# Copy the formatting already used by the other "已完成" rows (e.g. B2) so the
# green fill used throughout the list is reused instead of inventing a new one.
$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "已完成"

# Row 15 - 罐形容器 (tank container): Not Started -> In Progress (yellow fill).
$ws.Range("B15").Value = "进行中"
$ws.Range("B15").Interior.Color = 65535

$ws.Range("A1").Select()
